# "Generate Report for Archive"
#
# Localization status moved on from handoff into active translation, so:
#   1. every "Ready for handoff" status cell becomes "In Translation"
#      (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3)
#   2. the Status-ish columns that held the old, longer string re-fit to
#      the shorter text (Overview columns E/F, and column C on the two
#      per-locale detail sheets)

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: zh-cn / de-de status columns (E, F) -------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Narrower text -> narrower columns (re-fit E & F to the new content width).
$wsOverview.Range("E1:F3").ColumnWidth = 12.5

# --- zh-cn / de-de detail sheets: Status column (C) -------------------------
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus
    $ws.Range("C1:C3").ColumnWidth = 12.5
}
